$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: test_step_table (existing sheet, data model updated)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Clear old leftover columns (F:H) from the previous data model.
$ws1.Range("F1:H3").Clear()

$step1 = @("random-step-1-id", "random-test-id-1", "random-action-id-1", "random description",   "random data 1")
$step2 = @("random-step-2-id", "random-test-id-1", "random-action-id-2", "random description 2", "random data 2")
$step3 = @("random-step-3-id", "random-test-id-1", "random-action-id-3", "random description 3", "random data 3")
$stepRows = @($step1, $step2, $step3)

for ($r = 0; $r -lt $stepRows.Length; $r++) {
    $row = $stepRows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws1.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$ws1.Columns.Item(1).ColumnWidth = 14.59
$ws1.Columns.Item(3).ColumnWidth = 16.40
$ws1.Columns.Item(4).ColumnWidth = 18.00

# ---------------------------------------------------------------------------
# Sheet 2: test_case_table -> renamed to test_action_table, data model updated
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "test_action_table"

$action1 = @("random-action-id-1", "random action 1 name", "random action description 1", "random_action_runner_1_name")
$action2 = @("random-action-id-2", "random action 2 name", "random action description 2", "random_action_runner_2_name")
$action3 = @("random-action-id-3", "random action 3 name", "random action description 3", "random_action_runner_3_name")
$actionRows = @($action1, $action2, $action3)

for ($r = 0; $r -lt $actionRows.Length; $r++) {
    $row = $actionRows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws2.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$ws2.Columns.Item(1).ColumnWidth = 16.40
$ws2.Columns.Item(2).ColumnWidth = 17.60
$ws2.Columns.Item(3).ColumnWidth = 23.60
$ws2.Columns.Item(4).ColumnWidth = 27.50

# ---------------------------------------------------------------------------
# Sheet 3: test_case_table (new sheet, placed after test_action_table)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "test_case_table"

$ws3.Cells.Item(1, 1).Value = "random-test-id-1"
$ws3.Cells.Item(1, 2).Value = "description"
